$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "R" across columns C:K for rows 5, 6, and 7 (same value already used
# elsewhere in the sheet, e.g. K4, H27, I20, J20, etc.)
foreach ($r in 5..7) {
    $ws.Range("C$r" + ":K$r").Value = "R"
}

# Update the active selection to K4, matching the saved view state
$ws.Range("K4").Select()
